# Adds GUI elements / pop-up status information messaging system.
# - Typography sheet: change "Small" typography font to isocpeur.ttf and
#   add a new "XXLarge" typography row.
# - Translation sheet: fix a few existing Slovenian (SI) translations and
#   add a block of new TEXT IDs used by the new status-message popup and
#   other new GUI elements (Reset button, unit labels, language labels...).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Typography sheet
# ---------------------------------------------------------------------
$typo = $wb.Worksheets.Item("Typography")

# Row 6 ("Small") switches font from verdana.ttf to isocpeur.ttf
$typo.Cells.Item(6, 3).Value = "isocpeur.ttf"

# New row 8: "XXLarge" typography entry
$typo.Cells.Item(8, 2).Value = "XXLarge"
$typo.Cells.Item(8, 3).Value = "isocpeur.ttf"
$typo.Cells.Item(8, 4).Value = 60
$typo.Cells.Item(8, 5).Value = 4
$typo.Cells.Item(8, 6).Value = "?"
$typo.Cells.Item(8, 9).Value = "0-9"

# ---------------------------------------------------------------------
# Translation sheet
# ---------------------------------------------------------------------
$tr = $wb.Worksheets.Item("Translation")

# Fix a handful of existing Slovenian (SI / column G) translations.
# NOTE: a lone "0" looks numeric, so prefix with an apostrophe to force it
# to be stored as text (matching the other text cells in this column)
# instead of being auto-converted to a Number by Excel.
$tr.Cells.Item(29, 7).Value = "<value> kos"
$tr.Cells.Item(30, 7).Value = "'0"
$tr.Cells.Item(31, 7).Value = "<value> mm"
$tr.Cells.Item(32, 7).Value = "'0"
$tr.Cells.Item(33, 7).Value = "<value> mm"
$tr.Cells.Item(34, 7).Value = "'0"

# New rows 37-52: new GUI text ids (status messages, reset button, units...)
# Columns: B=TEXT ID, C=TYPOGRAPHY NAME, D=ALIGNMENT, E=DIRECTION, F=GB, G=SI
$newRows = @(
    @{ Row=37; B="STATUSMSG_DELAY_OF";           C="Default"; D="Center"; E="LTR"; F="Sum of delay and duration should not exceed 42 seconds. The delay was trimmed and relay will never be turned on.";          G="Vsota zamika in trajanja ne sme preseči 42 sekund. Zamik je bil obrezan in rele ne bo nikoli vklopljen." }
    @{ Row=38; B="STATUSMSG_DURATION_OF";        C="Default"; D="Center"; E="LTR"; F="Sum of delay and duration should not exceed 42 seconds. The duration was trimmed and relay will be turned off earlier.";     G="Vsota zamika in trajanja ne sme preseči 42 sekund. Trajanje je bilo omejeno in rele bo izklopljen predčasno." }
    @{ Row=39; B="STATUSMSG_RELAY_DEACTIVATED";  C="Default"; D="Center"; E="LTR"; F="Duration is 0. Relay is now inactive.";                                                                                       G="Trajanje je 0. Rele je neaktiven." }
    @{ Row=40; B="STATUSMSG_SET_LENGTH_TRIMMED"; C="Default"; D="Center"; E="LTR"; F="Maximum allowed length is 10000 mm. The length was trimmmed.";                                                                G="Najdaljsa dovoljena dolzina je 10000 mm. Dolzina je bila skrajsana." }
    @{ Row=41; B="STATUSMSG_OTHER_ERR";          C="Default"; D="Center"; E="LTR"; F="Unhandled error.";                                                                                                            G="Neznana napaka." }
    @{ Row=42; B="SingleUseId45";                C="Small";   D="Center"; E="LTR"; F="Status id: <value>";                                                                                                          G="Status id: <value>" }
    @{ Row=43; B="SingleUseId46";                C="Small";   D="Left";   E="LTR"; F="'0";                                                                                                                          G="-" }
    @{ Row=44; B="SingleUseId47";                C="Large";   D="Left";   E="LTR"; F="New Text";                                                                                                                    G="New Text" }
    @{ Row=45; B="SingleUseId48";                C="Large";   D="Left";   E="LTR"; F=".";                                                                                                                           G="." }
    @{ Row=46; B="SingleUseId49";                C="Large";   D="Center"; E="LTR"; F="SLO";                                                                                                                         G="SLO" }
    @{ Row=47; B="SingleUseId50";                C="Large";   D="Left";   E="LTR"; F="mm";                                                                                                                          G="mm" }
    @{ Row=48; B="SingleUseId51";                C="Large";   D="Left";   E="LTR"; F=".";                                                                                                                           G="." }
    @{ Row=49; B="SingleUseId52";                C="Large";   D="Center"; E="LTR"; F="ENG";                                                                                                                         G="ENG" }
    @{ Row=50; B="SingleUseId53";                C="XXLarge"; D="Center"; E="LTR"; F="<value>m";                                                                                                                    G="<value>m" }
    @{ Row=51; B="SingleUseId54";                C="XXLarge"; D="Left";   E="LTR"; F="'0";                                                                                                                          G="'0" }
    @{ Row=52; B="SingleUseId55";                C="Large";   D="Center"; E="LTR"; F="Reset";                                                                                                                       G="Ponastavi" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $tr.Cells.Item($row, 2).Value = $r.B
    $tr.Cells.Item($row, 3).Value = $r.C
    $tr.Cells.Item($row, 4).Value = $r.D
    $tr.Cells.Item($row, 5).Value = $r.E
    $tr.Cells.Item($row, 6).Value = $r.F
    $tr.Cells.Item($row, 7).Value = $r.G
}
